# Update "想去人数" (want-to-go count) values in column F for specific rows
# on both the "展览" sheet and the "全部类型" sheet, matching the upstream
# data refresh recorded in the commit "Update gh-pages to output generated
# at 456a3b4".

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# "展览" sheet (sheet1.xml): rows 3,5,7,11,12
$sheetExhibition.Range("F3").Value  = 996
$sheetExhibition.Range("F5").Value  = 2689
$sheetExhibition.Range("F7").Value  = 206
$sheetExhibition.Range("F11").Value = 2521
$sheetExhibition.Range("F12").Value = 616

# "全部类型" sheet (sheet4.xml): rows 4,6,8,13,14
$sheetAllTypes.Range("F4").Value  = 996
$sheetAllTypes.Range("F6").Value  = 2689
$sheetAllTypes.Range("F8").Value  = 206
$sheetAllTypes.Range("F13").Value = 2521
$sheetAllTypes.Range("F14").Value = 616
